# Auto-generated edit script: applies the numeric cell updates described by the commit diff.
# (Profit/loss recompute for several recipes across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets.)
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 7459.3145  # H17
$ws.Cells.Item(17, 10).Value = 1796.3529  # J17
$ws.Cells.Item(17, 12).Value = 5389.0587  # L17
$ws.Cells.Item(17, 14).Value = -5725.0587  # N17
$ws.Cells.Item(58, 8).Value = 1309.25  # H58
$ws.Cells.Item(58, 9).Value = 344.5  # I58
$ws.Cells.Item(58, 11).Value = 1033.5  # K58
$ws.Cells.Item(58, 13).Value = -883.5  # M58
$ws.Cells.Item(69, 8).Value = 4200  # H69
$ws.Cells.Item(69, 10).Value = 4200  # J69
$ws.Cells.Item(69, 12).Value = 12600  # L69
$ws.Cells.Item(69, 14).Value = -14348  # N69
$ws.Cells.Item(72, 8).Value = 4200  # H72
$ws.Cells.Item(72, 10).Value = 4200  # J72
$ws.Cells.Item(72, 12).Value = 37800  # L72
$ws.Cells.Item(72, 14).Value = -46536  # N72
$ws.Cells.Item(82, 8).Value = 7350  # H82
$ws.Cells.Item(82, 9).Value = 6575  # I82
$ws.Cells.Item(82, 10).Value = 7866.6665  # J82
$ws.Cells.Item(82, 11).Value = 19725  # K82
$ws.Cells.Item(82, 12).Value = 23599.9995  # L82
$ws.Cells.Item(82, 13).Value = -19319  # M82
$ws.Cells.Item(82, 14).Value = -24411.9995  # N82
$ws.Cells.Item(85, 8).Value = 7350  # H85
$ws.Cells.Item(85, 9).Value = 6575  # I85
$ws.Cells.Item(85, 10).Value = 7866.6665  # J85
$ws.Cells.Item(85, 11).Value = 19725  # K85
$ws.Cells.Item(85, 12).Value = 23599.9995  # L85
$ws.Cells.Item(85, 13).Value = -18321  # M85
$ws.Cells.Item(85, 14).Value = -26407.9995  # N85
$ws.Cells.Item(96, 8).Value = 749.875  # H96
$ws.Cells.Item(96, 9).Value = 599.8  # I96
$ws.Cells.Item(96, 11).Value = 1799.4  # K96
$ws.Cells.Item(96, 13).Value = -426.3999999999999  # M96
$ws.Cells.Item(129, 8).Value = 4718336.5  # H129
$ws.Cells.Item(129, 9).Value = 35715412  # I129
$ws.Cells.Item(129, 10).Value = 1390.2826  # J129
$ws.Cells.Item(129, 11).Value = 107146236  # K129
$ws.Cells.Item(129, 12).Value = 4170.8478  # L129
$ws.Cells.Item(129, 13).Value = -107141236  # M129
$ws.Cells.Item(129, 14).Value = -14170.8478  # N129
$ws.Cells.Item(132, 8).Value = 5265578  # H132
$ws.Cells.Item(132, 9).Value = 5884677  # I132
$ws.Cells.Item(132, 10).Value = 3235  # J132
$ws.Cells.Item(132, 11).Value = 17654031  # K132
$ws.Cells.Item(132, 12).Value = 9705  # L132
$ws.Cells.Item(132, 13).Value = -17651501  # M132
$ws.Cells.Item(132, 14).Value = -14765  # N132
$ws.Cells.Item(135, 8).Value = 1591.4783  # H135
$ws.Cells.Item(135, 9).Value = 1177.6666  # I135
$ws.Cells.Item(135, 11).Value = 10598.9994  # K135
$ws.Cells.Item(135, 13).Value = -8063.999400000001  # M135
$ws.Cells.Item(137, 8).Value = 3822.484  # H137
$ws.Cells.Item(137, 9).Value = 3837.9524  # I137
$ws.Cells.Item(137, 10).Value = 3790  # J137
$ws.Cells.Item(137, 11).Value = 11513.8572  # K137
$ws.Cells.Item(137, 12).Value = 11370  # L137
$ws.Cells.Item(137, 13).Value = -8963.8572  # M137
$ws.Cells.Item(137, 14).Value = -16470  # N137
$ws.Cells.Item(138, 8).Value = 2772.5325  # H138
$ws.Cells.Item(138, 9).Value = 1541.4791  # I138
$ws.Cells.Item(138, 10).Value = 4810.1377  # J138
$ws.Cells.Item(138, 11).Value = 4624.4373  # K138
$ws.Cells.Item(138, 12).Value = 14430.4131  # L138
$ws.Cells.Item(138, 13).Value = 515.5627000000004  # M138
$ws.Cells.Item(138, 14).Value = -24710.4131  # N138
$ws.Cells.Item(141, 8).Value = 536045.0600000001  # H141
$ws.Cells.Item(141, 9).Value = 1516.1666  # I141
$ws.Cells.Item(141, 11).Value = 4548.4998  # K141
$ws.Cells.Item(141, 13).Value = 631.5002000000004  # M141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3718.45  # H32
$ws.Cells.Item(32, 9).Value = 3262.6296  # I32
$ws.Cells.Item(32, 10).Value = 5661.684  # J32
$ws.Cells.Item(32, 11).Value = 3262.6296  # K32
$ws.Cells.Item(32, 12).Value = 5661.684  # L32
$ws.Cells.Item(32, 13).Value = -2975.6296  # M32
$ws.Cells.Item(32, 14).Value = -6235.684  # N32
$ws.Cells.Item(53, 8).Value = 45471.5  # H53
$ws.Cells.Item(53, 9).Value = 1800  # I53
$ws.Cells.Item(53, 10).Value = 60028.668  # J53
$ws.Cells.Item(53, 11).Value = 1800  # K53
$ws.Cells.Item(53, 12).Value = 60028.668  # L53
$ws.Cells.Item(53, 13).Value = -1118  # M53
$ws.Cells.Item(53, 14).Value = -61392.668  # N53
$ws.Cells.Item(61, 8).Value = 2636.4707  # H61
$ws.Cells.Item(61, 9).Value = 1202  # I61
$ws.Cells.Item(61, 10).Value = 4685.7144  # J61
$ws.Cells.Item(61, 11).Value = 1202  # K61
$ws.Cells.Item(61, 12).Value = 4685.7144  # L61
$ws.Cells.Item(61, 13).Value = -990  # M61
$ws.Cells.Item(61, 14).Value = -5109.7144  # N61
$ws.Cells.Item(74, 8).Value = 817.35  # H74
$ws.Cells.Item(74, 9).Value = 797.2105  # I74
$ws.Cells.Item(74, 10).Value = 1200  # J74
$ws.Cells.Item(74, 11).Value = 797.2105  # K74
$ws.Cells.Item(74, 12).Value = 1200  # L74
$ws.Cells.Item(74, 13).Value = 76.78949999999998  # M74
$ws.Cells.Item(74, 14).Value = -2948  # N74
$ws.Cells.Item(77, 8).Value = 817.35  # H77
$ws.Cells.Item(77, 9).Value = 797.2105  # I77
$ws.Cells.Item(77, 10).Value = 1200  # J77
$ws.Cells.Item(77, 11).Value = 3986.0525  # K77
$ws.Cells.Item(77, 12).Value = 6000  # L77
$ws.Cells.Item(77, 13).Value = 381.9474999999998  # M77
$ws.Cells.Item(77, 14).Value = -14736  # N77
$ws.Cells.Item(132, 8).Value = 17859486  # H132
$ws.Cells.Item(132, 9).Value = 25642578  # I132
$ws.Cells.Item(132, 10).Value = 4158.4116  # J132
$ws.Cells.Item(132, 11).Value = 76927734  # K132
$ws.Cells.Item(132, 12).Value = 12475.2348  # L132
$ws.Cells.Item(132, 13).Value = -76925204  # M132
$ws.Cells.Item(132, 14).Value = -17535.2348  # N132
$ws.Cells.Item(136, 8).Value = 2636.4707  # H136
$ws.Cells.Item(136, 9).Value = 1202  # I136
$ws.Cells.Item(136, 10).Value = 4685.7144  # J136
$ws.Cells.Item(136, 11).Value = 3606  # K136
$ws.Cells.Item(136, 12).Value = 14057.1432  # L136
$ws.Cells.Item(136, 13).Value = -1056  # M136
$ws.Cells.Item(136, 14).Value = -19157.1432  # N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(56, 8).Value = 0  # H56
$ws.Cells.Item(56, 10).Value = 0  # J56
$ws.Cells.Item(56, 12).Value = 0  # L56

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2760.7368  # H31
$ws.Cells.Item(31, 9).Value = 2122.5454  # I31
$ws.Cells.Item(31, 10).Value = 6972.8  # J31
$ws.Cells.Item(31, 11).Value = 2122.5454  # K31
$ws.Cells.Item(31, 12).Value = 6972.8  # L31
$ws.Cells.Item(31, 13).Value = -1827.5454  # M31
$ws.Cells.Item(31, 14).Value = -7562.8  # N31
$ws.Cells.Item(34, 8).Value = 2760.7368  # H34
$ws.Cells.Item(34, 9).Value = 2122.5454  # I34
$ws.Cells.Item(34, 10).Value = 6972.8  # J34
$ws.Cells.Item(34, 11).Value = 2122.5454  # K34
$ws.Cells.Item(34, 12).Value = 6972.8  # L34
$ws.Cells.Item(34, 13).Value = -1920.5454  # M34
$ws.Cells.Item(34, 14).Value = -7376.8  # N34
$ws.Cells.Item(58, 8).Value = 14288530  # H58
$ws.Cells.Item(58, 9).Value = 1607.5  # I58
$ws.Cells.Item(58, 10).Value = 29415860  # J58
$ws.Cells.Item(58, 11).Value = 1607.5  # K58
$ws.Cells.Item(58, 12).Value = 29415860  # L58
$ws.Cells.Item(58, 13).Value = -1404.5  # M58
$ws.Cells.Item(58, 14).Value = -29416266  # N58
$ws.Cells.Item(59, 8).Value = 8709  # H59
$ws.Cells.Item(59, 9).Value = 0  # I59
$ws.Cells.Item(59, 10).Value = 8709  # J59
$ws.Cells.Item(59, 11).Value = 0  # K59
$ws.Cells.Item(59, 12).Value = 8709  # L59
$ws.Cells.Item(59, 14).Value = -10999  # N59
$ws.Cells.Item(132, 8).Value = 3062.625  # H132
$ws.Cells.Item(132, 9).Value = 2342.2632  # I132
$ws.Cells.Item(132, 10).Value = 5800  # J132
$ws.Cells.Item(132, 11).Value = 7026.7896  # K132
$ws.Cells.Item(132, 12).Value = 17400  # L132
$ws.Cells.Item(132, 13).Value = -4496.7896  # M132
$ws.Cells.Item(132, 14).Value = -22460  # N132
$ws.Cells.Item(134, 8).Value = 2293.7144  # H134
$ws.Cells.Item(134, 9).Value = 968  # I134
$ws.Cells.Item(134, 10).Value = 4680  # J134
$ws.Cells.Item(134, 11).Value = 2904  # K134
$ws.Cells.Item(134, 12).Value = 14040  # L134
$ws.Cells.Item(134, 13).Value = -369  # M134
$ws.Cells.Item(134, 14).Value = -19110  # N134
$ws.Cells.Item(136, 8).Value = 14288530  # H136
$ws.Cells.Item(136, 9).Value = 1607.5  # I136
$ws.Cells.Item(136, 10).Value = 29415860  # J136
$ws.Cells.Item(136, 11).Value = 4822.5  # K136
$ws.Cells.Item(136, 12).Value = 88247580  # L136
$ws.Cells.Item(136, 13).Value = -2272.5  # M136
$ws.Cells.Item(136, 14).Value = -88252680  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 13006.923  # H87
$ws.Cells.Item(87, 10).Value = 15677.777  # J87
$ws.Cells.Item(87, 12).Value = 47033.331  # L87
$ws.Cells.Item(87, 14).Value = -49529.331  # N87
$ws.Cells.Item(90, 8).Value = 13006.923  # H90
$ws.Cells.Item(90, 10).Value = 15677.777  # J90
$ws.Cells.Item(90, 12).Value = 141099.993  # L90
$ws.Cells.Item(90, 14).Value = -153579.993  # N90
$ws.Cells.Item(120, 8).Value = 19558.889  # H120
$ws.Cells.Item(120, 9).Value = 18015  # I120
$ws.Cells.Item(120, 11).Value = 54045  # K120
$ws.Cells.Item(120, 13).Value = -49207  # M120
$ws.Cells.Item(129, 8).Value = 21475.482  # H129
$ws.Cells.Item(129, 9).Value = 5792  # I129
$ws.Cells.Item(129, 10).Value = 25039.908  # J129
$ws.Cells.Item(129, 11).Value = 17376  # K129
$ws.Cells.Item(129, 12).Value = 75119.724  # L129
$ws.Cells.Item(129, 13).Value = -12376  # M129
$ws.Cells.Item(129, 14).Value = -85119.724  # N129
$ws.Cells.Item(134, 8).Value = 2863.9167  # H134
$ws.Cells.Item(134, 9).Value = 1835.8  # I134
$ws.Cells.Item(134, 10).Value = 3598.2856  # J134
$ws.Cells.Item(134, 11).Value = 5507.4  # K134
$ws.Cells.Item(134, 12).Value = 10794.8568  # L134
$ws.Cells.Item(134, 13).Value = -437.3999999999996  # M134
$ws.Cells.Item(134, 14).Value = -20934.8568  # N134
$ws.Cells.Item(138, 8).Value = 1558  # H138
$ws.Cells.Item(138, 10).Value = 2604.6  # J138
$ws.Cells.Item(138, 12).Value = 7813.799999999999  # L138
$ws.Cells.Item(138, 14).Value = -18093.8  # N138

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2625  # H40
$ws.Cells.Item(40, 9).Value = 1000  # I40
$ws.Cells.Item(40, 10).Value = 3166.6667  # J40
$ws.Cells.Item(40, 11).Value = 1000  # K40
$ws.Cells.Item(40, 12).Value = 3166.6667  # L40
$ws.Cells.Item(40, 13).Value = -864  # M40
$ws.Cells.Item(40, 14).Value = -3438.6667  # N40
$ws.Cells.Item(46, 8).Value = 2172.5833  # H46
$ws.Cells.Item(46, 9).Value = 567.2857  # I46
$ws.Cells.Item(46, 10).Value = 4420  # J46
$ws.Cells.Item(46, 11).Value = 567.2857  # K46
$ws.Cells.Item(46, 12).Value = 4420  # L46
$ws.Cells.Item(46, 13).Value = -379.2857  # M46
$ws.Cells.Item(46, 14).Value = -4796  # N46
$ws.Cells.Item(68, 8).Value = 1447.1428  # H68
$ws.Cells.Item(68, 9).Value = 1019.5  # I68
$ws.Cells.Item(68, 10).Value = 10000  # J68
$ws.Cells.Item(68, 11).Value = 1019.5  # K68
$ws.Cells.Item(68, 12).Value = 10000  # L68
$ws.Cells.Item(68, 13).Value = -270.5  # M68
$ws.Cells.Item(68, 14).Value = -11498  # N68
$ws.Cells.Item(71, 8).Value = 1447.1428  # H71
$ws.Cells.Item(71, 9).Value = 1019.5  # I71
$ws.Cells.Item(71, 10).Value = 10000  # J71
$ws.Cells.Item(71, 11).Value = 5097.5  # K71
$ws.Cells.Item(71, 12).Value = 50000  # L71
$ws.Cells.Item(71, 13).Value = -1353.5  # M71
$ws.Cells.Item(71, 14).Value = -57488  # N71
$ws.Cells.Item(82, 8).Value = 2549.9  # H82
$ws.Cells.Item(82, 9).Value = 1722.1111  # I82
$ws.Cells.Item(82, 11).Value = 1722.1111  # K82
$ws.Cells.Item(82, 13).Value = -1361.1111  # M82
$ws.Cells.Item(85, 8).Value = 2549.9  # H85
$ws.Cells.Item(85, 9).Value = 1722.1111  # I85
$ws.Cells.Item(85, 11).Value = 1722.1111  # K85
$ws.Cells.Item(85, 13).Value = -474.1111000000001  # M85
$ws.Cells.Item(122, 8).Value = 3753.3333  # H122
$ws.Cells.Item(122, 9).Value = 2785.7144  # I122
$ws.Cells.Item(122, 10).Value = 4600  # J122
$ws.Cells.Item(122, 11).Value = 8357.143199999999  # K122
$ws.Cells.Item(122, 12).Value = 13800  # L122
$ws.Cells.Item(122, 13).Value = -5907.143199999999  # M122
$ws.Cells.Item(122, 14).Value = -18700  # N122
$ws.Cells.Item(136, 8).Value = 2072.8333  # H136
$ws.Cells.Item(136, 9).Value = 1506.6072  # I136
$ws.Cells.Item(136, 10).Value = 10000  # J136
$ws.Cells.Item(136, 11).Value = 4519.821599999999  # K136
$ws.Cells.Item(136, 12).Value = 30000  # L136
$ws.Cells.Item(136, 13).Value = -1969.821599999999  # M136
$ws.Cells.Item(136, 14).Value = -35100  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3554.138  # H132
$ws.Cells.Item(132, 10).Value = 12291.5  # J132
$ws.Cells.Item(132, 12).Value = 36874.5  # L132
$ws.Cells.Item(132, 14).Value = -41934.5  # N132

# Cells removed entirely by the target revision (no longer present in the row)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(56, 14).ClearContents()  # N56
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(59, 13).ClearContents()  # M59
